$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 659/660, pushing the existing rows 659:672 down to 661:674
$ws.Range("A659:A660").EntireRow.Insert()

# New row 659
$ws.Cells.Item(659, 1).Value = 9
$ws.Cells.Item(659, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(659, 3).Value = "Metropolitana"
$ws.Cells.Item(659, 4).Value = 44939
$ws.Cells.Item(659, 5).Value = 13
$ws.Cells.Item(659, 6).Value = 100112024
$ws.Cells.Item(659, 7).Value = "Choclo"
$ws.Cells.Item(659, 8).Value = "Choclero"
$ws.Cells.Item(659, 9).Value = "Primera"
$ws.Cells.Item(659, 10).Value = 10400
$ws.Cells.Item(659, 11).Value = 280
$ws.Cells.Item(659, 12).Value = 320
$ws.Cells.Item(659, 13).Value = 300
$ws.Cells.Item(659, 14).Value = "$/unidad"
$ws.Cells.Item(659, 15).Value = "Provincia de Melipilla"
$ws.Cells.Item(659, 16).Value = 300
$ws.Cells.Item(659, 17).Value = 1
$ws.Cells.Item(659, 18).Value = "Hortaliza"

# New row 660
$ws.Cells.Item(660, 1).Value = 9
$ws.Cells.Item(660, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(660, 3).Value = "Metropolitana"
$ws.Cells.Item(660, 4).Value = 44939
$ws.Cells.Item(660, 5).Value = 13
$ws.Cells.Item(660, 6).Value = 100112024
$ws.Cells.Item(660, 7).Value = "Choclo"
$ws.Cells.Item(660, 8).Value = "Choclero"
$ws.Cells.Item(660, 9).Value = "Primera"
$ws.Cells.Item(660, 10).Value = 9700
$ws.Cells.Item(660, 11).Value = 280
$ws.Cells.Item(660, 12).Value = 320
$ws.Cells.Item(660, 13).Value = 300
$ws.Cells.Item(660, 14).Value = "$/unidad"
$ws.Cells.Item(660, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(660, 16).Value = 300
$ws.Cells.Item(660, 17).Value = 1
$ws.Cells.Item(660, 18).Value = "Hortaliza"
